$wb = $excel.ActiveWorkbook

# 1) Rename the "Requested quantity" headers on the two existing sheets.
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 2) Add the new "PO Forecast" sheet (will land right after "Monthly Trend"
#    once we move it at the end, matching sheetId=3 / last tab position).
$wsForecast = $wb.Worksheets.Add()
$wsForecast.Name = "PO Forecast"

# Header row.
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$header = $wsForecast.Range("A1:D1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

# Data rows: ds (date-serial), PO_Forecast, yhat_lower, yhat_upper
$data = @(
    @(45137.99999999999,90,-26.77093327242614,197.1036203628788),
    @(45144.99999999999,91,-28.34435687463388,207.8012075928148),
    @(45158.99999999999,91,-25.00287120232229,215.0048090372464),
    @(45165.99999999999,92,-32.49815683967188,210.7653738802321),
    @(45179.99999999999,93,-14.92366384431766,209.4061595613706),
    @(45200.99999999999,94,-26.24933489924188,209.3711668940096),
    @(45207.99999999999,94,-21.86625020045173,218.5462916914325),
    @(45214.99999999999,94,-16.23502117144345,215.9529779981974),
    @(45221.99999999999,95,-18.17853770308538,215.4524976091776),
    @(45228.99999999999,95,-24.90960848530657,211.0861733964272),
    @(45235.99999999999,96,-27.82983061299649,213.7544628802802),
    @(45242.99999999999,96,-22.68703067029572,204.5180303995768),
    @(45249.99999999999,96,-23.57791828223142,214.655379359335),
    @(45256.99999999999,97,-21.94801670307535,213.4135791371891),
    @(45263.99999999999,97,-29.52334229425669,218.3971349698443),
    @(45270.99999999999,98,-23.23096555389341,209.9529126406403),
    @(45277.99999999999,98,-15.6058151912441,215.7791296110121),
    @(45298.99999999999,99,-13.93347969246602,228.1001671279696),
    @(45312.99999999999,100,-21.53004958963809,222.8758857606939),
    @(45319.99999999999,100,-17.67023145991554,223.7289205864072),
    @(45326.99999999999,101,-20.38312576890415,214.0526110220767),
    @(45333.99999999999,101,-8.106147118195246,225.0447793476864),
    @(45347.99999999999,102,-13.3536064660002,219.0930819703742),
    @(45354.99999999999,102,-13.54960859064898,218.095429679047),
    @(45361.99999999999,103,-10.66226555288068,222.9067643708891),
    @(45368.99999999999,103,-9.829866958872639,219.3249956532548),
    @(45375.99999999999,104,-5.419621121145923,225.2994468414389),
    @(45382.99999999999,104,-11.40050020784378,230.9408855329721),
    @(45389.99999999999,104,-0.07354437002106781,217.387668296385),
    @(45396.99999999999,105,-10.51953142041465,223.6116570083973),
    @(45403.99999999999,105,-19.30462578699729,220.3338567128682),
    @(45410.99999999999,105,1.564818926640104,232.2924919176567),
    @(45417.99999999999,106,-9.670169872691014,228.5306240731021),
    @(45424.99999999999,106,-9.780694146417416,211.6562361185798),
    @(45431.99999999999,107,-10.4458088584146,228.7779752623345),
    @(45438.99999999999,107,-17.25245505794212,226.7884921378349),
    @(45445.99999999999,107,-14.15573141290012,221.5198974328387),
    @(45452.99999999999,108,-12.2508455566622,227.7575086767751),
    @(45459.99999999999,108,-10.87212763030814,217.5562546763895),
    @(45466.99999999999,109,-13.40478912522432,232.5302084031013),
    @(45473.99999999999,109,-9.150040054599762,228.9890211337084)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Match the page margins used on the other sheets (0.75/0.75/1/1/0.5/0.5 in).
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# 3) Move the new sheet to the end of the workbook (after "Monthly Trend").
$lastExisting = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast.Move($null, $lastExisting)

Write-Output "PO Forecast sheet added with $($data.Count) data rows."
